# fixed filter issues + started SLURM implementation
#
# The "VAFRatio" filter column (L) is renamed/replaced by the already-present
# "VAFSim" filter semantics (matching the AMLValidation sheet which already
# used VAFSim/exclude similarVAF), the associated threshold values move from
# a "maximum ratio" style (1 / 0.9 / 0.8) to a "similarity" style (0 / 0.1 /
# 0.2), and the NHL sheet gains the PopFreq (O) column formatting/values that
# AMLMono7 and AMLValidation already had.

$wb = $excel.ActiveWorkbook

# --- AMLMono7 (sheet 1) -----------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("L1").Value = "VAFSim"
$ws1.Range("L3").Value = 0
$ws1.Range("L4").Value = 0.1
$ws1.Range("L5").Value = 0.2
$ws1.Range("L7").Value = "exclude similarVAF"

# --- NHL (sheet 2) ------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("L1").Value = "VAFSim"
$ws2.Range("L3").Value = 0
$ws2.Range("L4").Value = 0.1
$ws2.Range("L5").Value = 0.2
$ws2.Range("L7").Value = "exclude similarVAF"

# NHL didn't have the PopFreq (O) column filled in/formatted yet - bring it
# in line with the other two sheets.
$ws2.Range("O2").Font.Color = 0
$ws2.Range("O3").Font.Color = 0
$ws2.Range("O4").Font.Color = 0
$ws2.Range("O5").Value = 0.0001
$ws2.Range("O5").Font.Color = 0

# --- AMLValidation (sheet 3) --------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# widen the Tdepth-ish column so the longer header fits
$ws3.Columns.Item(6).ColumnWidth = 14.65

# --- window / selection bookkeeping -------------------------------------
$ws1.Activate()
$ws1.Range("K21").Select()

$ws2.Activate()
$ws2.Range("L1:L8").Select()

$ws3.Activate()
$ws3.Range("L1:L8").Select()

# leave the first sheet as the active / selected tab, as in the original file
$ws1.Activate()
